$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Append new row 229 as a full copy of row 228 (values + style),
# matching the workbook growing from 228 to 229 data rows.
$ws.Range("A228:R228").Copy($ws.Range("A229:R229"))

# Step 2: Shift the D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) values of rows 150-228 down
# by one row (row N receives what row N-1 originally held), since a new daily
# record is being inserted at the top of this block (row 150).
$dCol = @(44671,44336,44893,44782,44894,44680,44705,44179,44448,44529,44685,44221,44397,44908,44460,44319,44396,44729,44342,44532,44428,44767,44785,44322,44799,44344,44438,44634,44426,44441,44638,44270,44406,44907,44882,44663,44578,44490,44895,44895,44386,44425,44211,44452,44400,44776,44216,44376,44847,44839,44629,44299,44382,44756,44364,44453,44832,44162,44706,44883,44879,44665,44477,44813,44217,44623,44679,44414,44350,44763,44245,44809,44622,44757,44650,44357,44911,44736,44412)
$jCol = @(15,50,40,140,20,80,50,100,55,15,50,55,30,80,30,20,80,45,40,30,120,50,20,15,80,50,20,70,50,30,50,100,20,70,190,40,200,100,40,20,30,30,35,130,40,30,30,15,25,20,15,40,30,100,40,40,20,130,40,80,90,60,30,40,80,160,50,80,25,25,45,120,40,80,60,50,80,40,30)
$kCol = @(20000,28000,20000,16000,20000,20000,20000,33000,20000,20000,20000,35000,21000,20000,20000,35000,21000,31000,30000,20000,20000,20000,16000,35000,20000,25000,30000,22000,30000,25000,20000,40000,25000,20000,18000,25000,20000,25000,20000,20000,25000,25000,35000,20000,20000,20000,35000,25000,20000,20000,25000,35000,25000,20000,25000,20000,20000,33000,20000,20000,18000,15000,20000,20000,35000,22000,20000,20000,28000,20000,35000,20000,26000,20000,20000,35000,20000,20000,20000)
$lCol = @(20000,28000,20000,20000,20000,20000,20000,33000,20000,20000,20000,35000,21000,20000,20000,35000,21000,31000,30000,20000,25000,20000,16000,35000,20000,25000,30000,25000,30000,25000,20000,40000,25000,20000,20000,25000,20000,25000,20000,20000,25000,25000,35000,25000,20000,20000,35000,25000,20000,20000,25000,35000,25000,20000,25000,25000,20000,35000,20000,20000,20000,25000,20000,20000,35000,25000,20000,20000,28000,20000,35000,20000,26000,20000,25000,35000,20000,20000,20000)
$mCol = @(20000,28000,20000,17429,20000,20000,20000,33000,20000,20000,20000,35000,21000,20000,20000,35000,21000,31000,30000,20000,22083,20000,16000,35000,20000,25000,30000,23714,30000,25000,20000,40000,25000,20000,18947,25000,20000,25000,20000,20000,25000,25000,35000,23077,20000,20000,35000,25000,20000,20000,25000,35000,25000,20000,25000,22500,20000,33769,20000,20000,18889,21667,20000,20000,35000,23500,20000,20000,28000,20000,35000,20000,26000,20000,21667,35000,20000,20000,20000)
$pCol = @(1538,2154,1538,1341,1538,1538,1538,2538,1538,1538,1538,2692,1615,1538,1538,2692,1615,2385,2308,1538,1699,1538,1231,2692,1538,1923,2308,1824,2308,1923,1538,3077,1923,1538,1457,1923,1538,1923,1538,1538,1923,1923,2692,1775,1538,1538,2692,1923,1538,1538,1923,2692,1923,1538,1923,1731,1538,2598,1538,1538,1453,1667,1538,1538,2692,1808,1538,1538,2154,1538,2692,1538,2000,1538,1667,2692,1538,1538,1538)

# Walk from the bottom (row 228) up to row 151 so each row is overwritten with
# the previous row value before that source row itself gets overwritten.
for ($i = $dCol.Length - 1; $i -ge 1; $i--) {
    $targetRow = 150 + $i
    $ws.Range("D$targetRow").Value = $dCol[$i - 1]
    $ws.Range("J$targetRow").Value = $jCol[$i - 1]
    $ws.Range("K$targetRow").Value = $kCol[$i - 1]
    $ws.Range("L$targetRow").Value = $lCol[$i - 1]
    $ws.Range("M$targetRow").Value = $mCol[$i - 1]
    $ws.Range("P$targetRow").Value = $pCol[$i - 1]
}

# Step 3: Row 150 becomes the new, most-recent daily record (new Fecha/Volumen);
# K, L, M, P are unchanged for this row.
$ws.Range("D150").Value = 44917
$ws.Range("J150").Value = 35
